# Add a new "Swiss" test-data sheet, modeled on the existing market sheets.
#
# The new sheet needs:
#   - the same column widths as Germany/Belgium (B=38.88671875, C=13.21875, D=22.77734375)
#   - the same 12-row "Loops" layout as Germany/Czech (XLM800 / XLM800-STI /
#     XLM800-Zetfas / Wg / Loops)
#   - B2 = "Switzerland Market", B4 = "NGC-3476/T2653"
#
# Duplicating Belgium preserves the correct column widths exactly (copying
# Czech and then rewriting ColumnWidth loses precision because ColumnWidth is
# character-based and rounds), so we copy Belgium and then insert the two
# extra "Loops" rows Belgium is missing, copying formatting down from the row
# above.

$wb = $excel.ActiveWorkbook

$belgium = $wb.Worksheets.Item("Belgium")
$belgium.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"

# Belgium's loop block is only 4 rows (Loops/XLM800/Wg/Loops); grow it to
# match the 6-row block used by Germany/Czech by inserting two rows for
# XLM800-STI / XLM800-Zetfas above the existing "Wg" row.
$swiss.Rows.Item(9).Resize(2).Insert()
$swiss.Range("A9").Value = "XLM800-STI"
$swiss.Range("A10").Value = "XLM800-Zetfas"

# Copy the bordered-cell formatting from the row above down onto the two new
# rows (Insert doesn't carry formatting along by itself).
$swiss.Range("A8").Copy()
$swiss.Range("A9:A10").PasteSpecial(-4122)

# Market-specific content.
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2653"

# Make this the active sheet/selection, matching the committed file.
$swiss.Activate()
$swiss.Range("A9").Select()

# The Czech sheet's selection was left as "select all" (e.g. Ctrl+A) in the
# committed workbook.
$czech = $wb.Worksheets.Item("Czech")
$czech.Activate()
$czech.Range("A1:XFD1048576").Select()

# Swiss is the sheet that ends up active/selected in the saved workbook.
$swiss.Activate()
